# Applies scheduled-runner value updates to the Leve profit/price columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1406.2
$ws.Range("I28").Value = 1406.2
$ws.Range("K28").Value = 1406.2
$ws.Range("M28").Value = -921.2

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H92").Value = 48178.094
$ws.Range("I92").Value = 53164
$ws.Range("J92").Value = 812
$ws.Range("K92").Value = 53164
$ws.Range("L92").Value = 812
$ws.Range("M92").Value = -51916
$ws.Range("N92").Value = -3308

$ws.Range("H132").Value = 3795.85
$ws.Range("I132").Value = 3662.0557
$ws.Range("K132").Value = 10986.1671
$ws.Range("M132").Value = -8456.167099999999

$ws.Range("H141").Value = 15043
$ws.Range("I141").Value = 16673.375
$ws.Range("K141").Value = 50020.125
$ws.Range("M141").Value = -44840.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1302.6857
$ws.Range("I32").Value = 1302.6857
$ws.Range("K32").Value = 1302.6857
$ws.Range("M32").Value = -1015.6857

$ws.Range("H61").Value = 4645.054
$ws.Range("I61").Value = 3576.5356
$ws.Range("K61").Value = 3576.5356
$ws.Range("M61").Value = -3364.5356

$ws.Range("H132").Value = 4270.3335
$ws.Range("I132").Value = 4270.3335
$ws.Range("K132").Value = 12811.0005
$ws.Range("M132").Value = -10281.0005

$ws.Range("H136").Value = 4645.054
$ws.Range("I136").Value = 3576.5356
$ws.Range("K136").Value = 10729.6068
$ws.Range("M136").Value = -8179.606800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 19685
$ws.Range("I7").Value = 1000
$ws.Range("J7").Value = 29027.5
$ws.Range("K7").Value = 1000
$ws.Range("L7").Value = 29027.5
$ws.Range("M7").Value = -887
$ws.Range("N7").Value = -29253.5

$ws.Range("H20").Value = 1943.5
$ws.Range("I20").Value = 2076
$ws.Range("J20").Value = 1744.75
$ws.Range("K20").Value = 2076
$ws.Range("L20").Value = 1744.75
$ws.Range("M20").Value = -1829
$ws.Range("N20").Value = -2238.75

$ws.Range("H64").Value = 2796.8
$ws.Range("I64").Value = 1570.5
$ws.Range("J64").Value = 3103.375
$ws.Range("K64").Value = 1570.5
$ws.Range("L64").Value = 3103.375
$ws.Range("M64").Value = -1345.5
$ws.Range("N64").Value = -3553.375

$ws.Range("H67").Value = 2796.8
$ws.Range("I67").Value = 1570.5
$ws.Range("J67").Value = 3103.375
$ws.Range("K67").Value = 1570.5
$ws.Range("L67").Value = 3103.375
$ws.Range("M67").Value = -790.5
$ws.Range("N67").Value = -4663.375

$ws.Range("H86").Value = 8233.177
$ws.Range("I86").Value = 2413.9167
$ws.Range("K86").Value = 2413.9167
$ws.Range("M86").Value = -1290.9167

$ws.Range("H89").Value = 8233.177
$ws.Range("I89").Value = 2413.9167
$ws.Range("K89").Value = 12069.5835
$ws.Range("M89").Value = -6453.583500000001

$ws.Range("H95").Value = 15922.6
$ws.Range("J95").Value = 15922.6
$ws.Range("L95").Value = 15922.6
$ws.Range("N95").Value = -21414.6

$ws.Range("H105").Value = 3876
$ws.Range("I105").Value = 4401.8
$ws.Range("K105").Value = 4401.8
$ws.Range("M105").Value = -2654.8

$ws.Range("H134").Value = 9631.817999999999
$ws.Range("I134").Value = 10017.5
$ws.Range("K134").Value = 30052.5
$ws.Range("M134").Value = -27517.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4195.6665
$ws.Range("I58").Value = 2680.1428
$ws.Range("K58").Value = 2680.1428
$ws.Range("M58").Value = -2477.1428

$ws.Range("H136").Value = 4195.6665
$ws.Range("I136").Value = 2680.1428
$ws.Range("K136").Value = 8040.428400000001
$ws.Range("M136").Value = -5490.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1319.2
$ws.Range("I5").Value = 1386.1428
$ws.Range("K5").Value = 4158.428400000001
$ws.Range("M5").Value = -4046.428400000001

$ws.Range("H68").Value = 300
$ws.Range("I68").Value = 300
$ws.Range("K68").Value = 900
$ws.Range("M68").Value = -89

$ws.Range("H71").Value = 300
$ws.Range("I71").Value = 300
$ws.Range("K71").Value = 2700
$ws.Range("M71").Value = 1356

$ws.Range("H122").Value = 670.5714
$ws.Range("I122").Value = 585.125
$ws.Range("J122").Value = 784.5
$ws.Range("K122").Value = 5266.125
$ws.Range("L122").Value = 7060.5
$ws.Range("M122").Value = -2816.125
$ws.Range("N122").Value = -11960.5

$ws.Range("H135").Value = 1319.2
$ws.Range("I135").Value = 1386.1428
$ws.Range("K135").Value = 12475.2852
$ws.Range("M135").Value = -9940.2852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 23281.45
$ws.Range("I46").Value = 13961
$ws.Range("J46").Value = 25611.562
$ws.Range("K46").Value = 13961
$ws.Range("L46").Value = 25611.562
$ws.Range("M46").Value = -13805
$ws.Range("N46").Value = -25923.562

$ws.Range("H80").Value = 2504.5625
$ws.Range("I80").Value = 2520.7778
$ws.Range("J80").Value = 2483.7144
$ws.Range("K80").Value = 2520.7778
$ws.Range("L80").Value = 2483.7144
$ws.Range("M80").Value = -1522.7778
$ws.Range("N80").Value = -4479.7144

$ws.Range("H83").Value = 2504.5625
$ws.Range("I83").Value = 2520.7778
$ws.Range("J83").Value = 2483.7144
$ws.Range("K83").Value = 12603.889
$ws.Range("L83").Value = 12418.572
$ws.Range("M83").Value = -7611.888999999999
$ws.Range("N83").Value = -22402.572

$ws.Range("H122").Value = 3194.923
$ws.Range("I122").Value = 2836.1667
$ws.Range("K122").Value = 8508.500100000001
$ws.Range("M122").Value = -6058.500100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3372.2083
$ws.Range("I7").Value = 3024.2727
$ws.Range("J7").Value = 7199.5
$ws.Range("K7").Value = 3024.2727
$ws.Range("L7").Value = 7199.5
$ws.Range("M7").Value = -2912.2727
$ws.Range("N7").Value = -7423.5

$ws.Range("H126").Value = 3372.2083
$ws.Range("I126").Value = 3024.2727
$ws.Range("J126").Value = 7199.5
$ws.Range("K126").Value = 9072.8181
$ws.Range("L126").Value = 21598.5
$ws.Range("M126").Value = -6602.8181
$ws.Range("N126").Value = -26538.5

$ws.Range("H132").Value = 2041.4706
$ws.Range("J132").Value = 1899.75
$ws.Range("L132").Value = 5699.25
$ws.Range("N132").Value = -10759.25

$ws.Range("H136").Value = 6337.091
$ws.Range("I136").Value = 3467.5557
$ws.Range("K136").Value = 10402.6671
$ws.Range("M136").Value = -7852.667099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 47886.383
$ws.Range("I45").Value = 37343.168
$ws.Range("J45").Value = 56923.43
$ws.Range("K45").Value = 37343.168
$ws.Range("L45").Value = 56923.43
$ws.Range("M45").Value = -36852.168
$ws.Range("N45").Value = -57905.43

$ws.Range("H126").Value = 1169.3158
$ws.Range("I126").Value = 1096.1177
$ws.Range("J126").Value = 1791.5
$ws.Range("K126").Value = 3288.3531
$ws.Range("L126").Value = 5374.5
$ws.Range("M126").Value = -818.3531000000003
$ws.Range("N126").Value = -10314.5
